$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price/Volume columns to text so that numeric-looking
# strings (e.g. "0.110", "1.03") are not auto-converted to numbers by Excel,
# which would strip meaningful trailing zeros / alter formatting.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = '45.808.90'
$ws.Range("E2").Value = '  +7.42%  '
$ws.Range("D3").Value = '2.416.66'
$ws.Range("E3").Value = '  +5.11%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '116.32'
$ws.Range("E5").Value = '  +11.29%  '
$ws.Range("D6").Value = '319.37'
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("D7").Value = '0.641'
$ws.Range("E7").Value = '  +2.67%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +4.59%  '
$ws.Range("D10").Value = '43.27'
$ws.Range("E10").Value = '  +8.71%  '
$ws.Range("D11").Value = '0.0939'
$ws.Range("E11").Value = '  +4.56%  '
$ws.Range("D12").Value = '8.78'
$ws.Range("E12").Value = '  +6.77%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.110'
$ws.Range("E13").Value = '  +2.42%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '1.03'
$ws.Range("E14").Value = '  +4.10%  '
$ws.Range("D15").Value = '15.96'
$ws.Range("E15").Value = '  +4.07%  '
$ws.Range("D16").Value = '2.785.62'
$ws.Range("E16").Value = '  +5.27%  '
$ws.Range("D17").Value = '2.417.24'
$ws.Range("E17").Value = '  +4.97%  '
$ws.Range("D18").Value = '45.838.73'
$ws.Range("E18").Value = '  +7.55%  '
$ws.Range("D19").Value = '7.63'
$ws.Range("E19").Value = '  +4.25%  '
$ws.Range("E20").Value = '  +4.56%  '
$ws.Range("D21").Value = '13.48'
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").Value = '75.26'
$ws.Range("E22").Value = '  +2.51%  '
$ws.Range("D23").Value = '3.58'
$ws.Range("E23").Value = '  +4.36%  '
$ws.Range("D24").Value = '268.27'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '2.42'
$ws.Range("E25").Value = '  +10.32%  '
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '7.65'
$ws.Range("E27").Value = '  +7.10%  '
$ws.Range("D28").Value = '11.40'
$ws.Range("E28").Value = '  +5.40%  '
$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  +4.58%  '
$ws.Range("D30").Value = '40.23'
$ws.Range("E30").Value = '  +11.30%  '
$ws.Range("D31").Value = '23.05'
$ws.Range("E31").Value = '  +3.12%  '
$ws.Range("D32").Value = '0.0974'
$ws.Range("E32").Value = '  +14.15%  '
$ws.Range("D33").Value = '173.56'
$ws.Range("E33").Value = '  +5.39%  '
$ws.Range("D34").Value = '2.98'
$ws.Range("E34").Value = '  +13.29%  '
$ws.Range("D35").Value = '5.05'
$ws.Range("E35").Value = '  +10.81%  '
$ws.Range("D36").Value = '0.133'
$ws.Range("E36").Value = '  +2.63%  '
$ws.Range("E37").Value = '  +7.61%  '
$ws.Range("D38").Value = '4.29'
$ws.Range("E38").Value = '  +17.37%  '
$ws.Range("D39").Value = '3.15'
$ws.Range("E39").Value = '  +11.62%  '
$ws.Range("D40").Value = '0.0367'
$ws.Range("E40").Value = '  +5.96%  '
$ws.Range("D41").Value = '1.80'
$ws.Range("E41").Value = '  +14.32%  '
$ws.Range("D42").Value = '101.83'
$ws.Range("E42").Value = '  -5.56%  '
$ws.Range("D43").Value = '13.68'
$ws.Range("E43").Value = '  +12.58%  '
$ws.Range("D44").Value = '0.241'
$ws.Range("E44").Value = '  +6.03%  '
$ws.Range("D45").Value = '72.91'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '5.87'
$ws.Range("E47").Value = '  +14.09%  '
$ws.Range("D48").Value = '117.93'
$ws.Range("E48").Value = '  +6.55%  '
$ws.Range("D49").Value = '1.69'
$ws.Range("E49").Value = '  +18.66%  '
$ws.Range("D50").Value = '83.09'
$ws.Range("E50").Value = '  +6.87%  '
$ws.Range("D51").Value = '9.53'
$ws.Range("E51").Value = '  +10.15%  '

# Restore the original (default) style now that the text values are set,
# so no stray number formatting is left behind on the cells.
$priceVolumeRange.Style = "Normal"
